$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ConcentrationSingleValues")

$data = @(
    @("F01", "CPF", 0.02,   "mgPerKg"),
    @("F02", "CPF", 0.015,  "mgPerKg"),
    @("F03", "CPF", 0.03,   "mgPerKg"),
    @("F04", "CPF", 0.02,   "mgPerKg"),
    @("F05", "CPF", 0.01,   "mgPerKg"),
    @("F06", "CPF", 0.015,  "mgPerKg"),
    @("F07", "CPF", 0.03,   "mgPerKg"),
    @("F08", "CPF", 0.005,  "mgPerKg"),
    @("F09", "CPF", 0.002,  "mgPerKg"),
    @("F10", "CPF", 0.001,  "mgPerKg"),
    @("F11", "CPF", 0.001,  "mgPerKg"),
    @("F12", "CPF", 0.0005, "mgPerKg"),
    @("F13", "CPF", 0.0005, "mgPerKg"),
    @("F14", "CPF", 0.001,  "mgPerKg"),
    @("F15", "CPF", 0.001,  "mgPerKg"),
    @("F16", "CPF", 0.0008, "mgPerKg"),
    @("F17", "CPF", 0.001,  "mgPerKg"),
    @("F18", "CPF", 0.005,  "mgPerKg"),
    @("F19", "CPF", 0.015,  "mgPerKg"),
    @("F20", "CPF", 0.01,   "mgPerKg")
)

$reference = "EU 2022 monitoring (synthetic)"
$valueType = "MeanConcentration"

$row = 2
foreach ($item in $data) {
    $ws.Cells.Item($row, 1).Value = $item[0]
    $ws.Cells.Item($row, 2).Value = $item[1]
    $ws.Cells.Item($row, 3).Value = $item[2]
    $ws.Cells.Item($row, 4).Value = $valueType
    $ws.Cells.Item($row, 6).Value = $item[3]
    $ws.Cells.Item($row, 7).Value = $reference
    $row++
}

# Column E (Percentile) is left blank for every data row, but the source
# workbook still records an (empty) cell there, so force the cell to
# persist in the saved sheet by touching a per-cell property instead of
# a value.
$ws.Range("E2:E21").NumberFormat = "General"
